# "torch mle minor updates"
#
# Adds a normalization step to the nested-logit toy model on the
# "regular_simulation" sheet: z_gi^(1+theta) values (C49:E56) are now
# divided by a new normalizing constant in $H$61 before being raised to
# the power, to avoid infinity values in the probability computation.
# A small "testing" block (O79:Q100) is added to verify the ratios used
# downstream are unchanged by the normalization.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("regular_simulation")

# --- snapshot the pre-change probability ratios (I93:K100) -----------------
# These are the values the new "diff" check columns (O:Q) compare against;
# they must be captured before the H61 normalization changes I93:K100
# downstream (even though, being ratios, they end up identical).
$snapshot = @{}
for ($r = 93; $r -le 100; $r++) {
    $snapshot[$r] = @($ws.Range("I$r").Value2, $ws.Range("J$r").Value2, $ws.Range("K$r").Value2)
}

# --- new normalization constant --------------------------------------------
$ws.Range("H60").Value = "Testing for Normalization (to avoid infinity values in the probabilty computation)"
$ws.Range("H61").Value = 10

# --- rewrite C49:E56 to normalize by H61 before exponentiating --------------
for ($r = 49; $r -le 56; $r++) {
    $src = $r - 47
    foreach ($col in @("C", "D", "E")) {
        $ws.Range("$col$r").Formula = "=($col$src/`$H`$61)^(`$H`$48+1)"
    }
}

# --- "diff" testing block ----------------------------------------------------
$ws.Range("O79").Value = "diff"

for ($r = 80; $r -le 88; $r++) {
    $src = $r + 13
    $ws.Range("O$r").Formula = "=I$src-O$src"
    $ws.Range("P$r").Formula = "=J$src-P$src"
    $ws.Range("Q$r").Formula = "=K$src-Q$src"
}

$ws.Range("O92").Value = $ws.Range("I92").Value2

for ($r = 93; $r -le 100; $r++) {
    $vals = $snapshot[$r]
    $ws.Range("O$r").Value = $vals[0]
    $ws.Range("P$r").Value = $vals[1]
    $ws.Range("Q$r").Value = $vals[2]
    $ws.Range("O$r`:Q$r").NumberFormat = $ws.Range("I$r").NumberFormat
}
